$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain plain text so values like
# "1.00" or "0.0000132" keep their original formatting instead of being
# auto-converted to numbers by Excel's smart-entry parsing.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "54.220.59"
$ws.Range("E2").Value = "  -8.08%  "

# Row 3
$ws.Range("D3").Value = "2.412.18"
$ws.Range("E3").Value = "  -11.63%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").Value = "463.66"
$ws.Range("E5").Value = "  -8.32%  "

# Row 6
$ws.Range("D6").Value = "131.25"
$ws.Range("E6").Value = "  -7.64%  "

# Row 7
$ws.Range("E7").Value = "  -0.15%  "

# Row 8
$ws.Range("D8").Value = "0.489"
$ws.Range("E8").Value = "  -8.11%  "

# Row 9
$ws.Range("D9").Value = "2.429.25"
$ws.Range("E9").Value = "  -11.48%  "

# Row 10
$ws.Range("D10").Value = "0.0944"
$ws.Range("E10").Value = "  -9.95%  "

# Row 11
$ws.Range("D11").Value = "5.31"
$ws.Range("E11").Value = "  -12.76%  "

# Row 12
$ws.Range("D12").Value = "0.312"
$ws.Range("E12").Value = "  -10.32%  "

# Row 13
$ws.Range("E13").Value = "  -4.15%  "

# Row 14
$ws.Range("D14").Value = "2.842.95"
$ws.Range("E14").Value = "  -11.34%  "

# Row 15
$ws.Range("D15").Value = "54.093.40"
$ws.Range("E15").Value = "  -8.34%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000132"
$ws.Range("E16").Value = "  -3.32%  "

# Row 17
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "19.63"
$ws.Range("E17").Value = "  -9.64%  "

# Row 18
$ws.Range("D18").Value = "2.438.08"
$ws.Range("E18").Value = "  -10.68%  "

# Row 19
$ws.Range("D19").Value = "4.16"
$ws.Range("E19").Value = "  -12.85%  "

# Row 20
$ws.Range("D20").Value = "307.21"
$ws.Range("E20").Value = "  -10.73%  "

# Row 21
$ws.Range("D21").Value = "9.38"
$ws.Range("E21").Value = "  -14.93%  "

# Row 22
$ws.Range("D22").Value = "0.993"
$ws.Range("E22").Value = "  -0.32%  "

# Row 23
$ws.Range("E23").Value = "  +0.97%  "

# Row 24
$ws.Range("D24").Value = "5.32"
$ws.Range("E24").Value = "  -15.17%  "

# Row 25
$ws.Range("D25").Value = "55.92"
$ws.Range("E25").Value = "  -11.63%  "

# Row 26
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +1.09%  "

# Row 27
$ws.Range("D27").Value = "0.383"
$ws.Range("E27").Value = "  -10.47%  "

# Row 28
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.155"
$ws.Range("E28").Value = "  -9.82%  "

# Row 29
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.511.28"
$ws.Range("E29").Value = "  -11.92%  "

# Row 30
$ws.Range("D30").Value = "7.02"
$ws.Range("E30").Value = "  -7.04%  "

# Row 31
$ws.Range("E31").Value = "  -0.07%  "

# Row 32
$ws.Range("D32").Value = "0.0₃0710"
$ws.Range("E32").Value = "  -14.84%  "

# Row 33
$ws.Range("D33").Value = "146.11"
$ws.Range("E33").Value = "  -3.37%  "

# Row 34
$ws.Range("D34").Value = "17.65"
$ws.Range("E34").Value = "  -7.97%  "

# Row 35
$ws.Range("D35").Value = "1.43"
$ws.Range("E35").Value = "  -11.22%  "

# Row 36
$ws.Range("D36").Value = "4.96"
$ws.Range("E36").Value = "  -8.70%  "

# Row 37
$ws.Range("D37").Value = "3.51"
$ws.Range("E37").Value = "  -16.28%  "

# Row 38
$ws.Range("E38").Value = "  -7.38%  "

# Row 39
$ws.Range("D39").Value = "0.799"
$ws.Range("E39").Value = "  -16.28%  "

# Row 40
$ws.Range("E40").Value = "  -0.37%  "

# Row 41
$ws.Range("D41").Value = "32.85"
$ws.Range("E41").Value = "  -8.47%  "

# Row 42
$ws.Range("D42").Value = "0.593"
$ws.Range("E42").Value = "  -2.04%  "

# Row 43
$ws.Range("D43").Value = "0.0520"
$ws.Range("E43").Value = "  -7.24%  "

# Row 44
$ws.Range("D44").Value = "3.23"
$ws.Range("E44").Value = "  -9.31%  "

# Row 45
$ws.Range("D45").Value = "10.09"
$ws.Range("E45").Value = "  -2.71%  "

# Row 46
$ws.Range("D46").Value = "1.23"
$ws.Range("E46").Value = "  -12.07%  "

# Row 47
$ws.Range("D47").Value = "1.934.09"
$ws.Range("E47").Value = "  -11.74%  "

# Row 48
$ws.Range("D48").Value = "0.0871"
$ws.Range("E48").Value = "  -1.93%  "

# Row 49
$ws.Range("E49").Value = "  -4.64%  "

# Row 50
$ws.Range("E50").Value = "  -12.18%  "

# Row 51
$ws.Range("D51").Value = "16.49"
$ws.Range("E51").Value = "  -13.56%  "
